$d = $word.ActiveDocument

# The paragraph "--Circle around him waiting for an opportunity" is split
# across two runs: the main run, and a trailing " to perform a quick jab"
# run (separated by a _GoBack bookmark). We need to:
#   1. Update the first run's text so it reads "...opportunity to perform a
#      quick jab" (absorbing the trailing run's wording).
#   2. Remove the now-redundant trailing " to perform a quick jab" run.
#   3. Remove the whole following paragraph ("-- Do nothing and let him make
#      another move"), including its paragraph mark.

$circlePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("--Circle around him waiting for an opportunity")) {
        $circlePara = $p
        break
    }
}

$trailingText = " to perform a quick jab"
$paraEnd = $circlePara.Range.End
$trailStart = $paraEnd - 1 - $trailingText.Length
$trailRange = $d.Range($trailStart, $paraEnd - 1)
if ($trailRange.Text -eq $trailingText) {
    [void]$trailRange.Delete()
}

$circleRange = $circlePara.Range
[void]$circleRange.Find.Execute("--Circle around him waiting for an opportunity", $true, $false, $false, $false, $false, $true, 1, $false, "--Circle around him waiting for an opportunity to perform a quick jab", 2)

$doNothingPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("-- Do nothing and let him make another move")) {
        $doNothingPara = $p
        break
    }
}
[void]$doNothingPara.Range.Delete()
